# Update numeric values in column F ("弹幕数" / danmaku-like counter) on the
# "展览" and "全部类型" worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 501
$ws1.Range("F4").Value  = 422
$ws1.Range("F5").Value  = 8551
$ws1.Range("F7").Value  = 10727
$ws1.Range("F22").Value = 1820
$ws1.Range("F23").Value = 79
$ws1.Range("F25").Value = 344
$ws1.Range("F26").Value = 286
$ws1.Range("F31").Value = 22
$ws1.Range("F32").Value = 41
$ws1.Range("F33").Value = 1416
$ws1.Range("F34").Value = 444
$ws1.Range("F36").Value = 288
$ws1.Range("F37").Value = 23
$ws1.Range("F39").Value = 515
$ws1.Range("F43").Value = 642
$ws1.Range("F46").Value = 101

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 501
$ws4.Range("F8").Value  = 422
$ws4.Range("F9").Value  = 8551
$ws4.Range("F11").Value = 10727
$ws4.Range("F19").Value = 1820
$ws4.Range("F20").Value = 79
$ws4.Range("F22").Value = 344
$ws4.Range("F23").Value = 286
$ws4.Range("F29").Value = 22
$ws4.Range("F31").Value = 41
$ws4.Range("F34").Value = 1416
$ws4.Range("F35").Value = 444
$ws4.Range("F39").Value = 515
$ws4.Range("F46").Value = 642
$ws4.Range("F49").Value = 101

$wb.Save()
